$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.699.24"
$ws.Range("E2").Value = "  -1.31%  "
$ws.Range("D3").Value = "2.586.47"
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "551.96"
$ws.Range("E5").Value = "  +2.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.41"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.594"
$ws.Range("E8").Value = "  +3.93%  "
$ws.Range("E9").Value = "  +3.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.101"
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.140"
$ws.Range("E11").Value = "  +3.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.333"
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("D13").Value = "3.041.96"
$ws.Range("E13").Value = "  -2.26%  "
$ws.Range("D14").Value = "58.632.72"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.75"
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("D16").Value = "2.592.76"
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("E17").Value = "  -2.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.44"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "335.85"
$ws.Range("E19").Value = "  -1.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.03"
$ws.Range("E20").Value = "  -3.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.10"
$ws.Range("E21").Value = "  -2.81%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.50"
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.424"
$ws.Range("E24").Value = "  +2.30%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  -4.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.06"
$ws.Range("E27").Value = "  -2.98%  "
$ws.Range("D28").Value = "0.0₃0746"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +0.93%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.89"
$ws.Range("E31").Value = "  +0.99%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "154.18"
$ws.Range("E32").Value = "  +2.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.84"
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.90"
$ws.Range("E34").Value = "  -2.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "37.22"
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.855"
$ws.Range("E36").Value = "  +2.43%  "
$ws.Range("E37").Value = "  -2.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.44"
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.822"
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.59"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "279.61"
$ws.Range("E41").Value = "  -3.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.595"
$ws.Range("E43").Value = "  -1.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.65"
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0948"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0529"
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("D48").Value = "1.908.20"
$ws.Range("E48").Value = "  -3.25%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.82"
$ws.Range("E49").Value = "  -2.58%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.39"
$ws.Range("E50").Value = "  -3.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.71"
$ws.Range("E51").Value = "  +1.68%  "
